$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1965.2059
$ws.Range("I15").Value = 1965.2059
$ws.Range("K15").Value = 5895.6177
$ws.Range("M15").Value = -5726.6177
# row 28
$ws.Range("H28").Value = 556.8570999999999
$ws.Range("I28").Value = 556.8570999999999
$ws.Range("K28").Value = 556.8570999999999
$ws.Range("M28").Value = -71.85709999999995
# row 40
$ws.Range("H40").Value = 848.5
$ws.Range("I40").Value = 844
$ws.Range("K40").Value = 844
$ws.Range("M40").Value = -669
# row 64
$ws.Range("H64").Value = 1700
$ws.Range("J64").Value = 1700
$ws.Range("L64").Value = 1700
$ws.Range("N64").Value = -2196
# row 67
$ws.Range("H67").Value = 1700
$ws.Range("J67").Value = 1700
$ws.Range("L67").Value = 1700
$ws.Range("N67").Value = -3416
# row 113
$ws.Range("H113").Value = 18334666
$ws.Range("I113").Value = 2001599
$ws.Range("J113").Value = 100000000
$ws.Range("K113").Value = 2001599
$ws.Range("L113").Value = 100000000
$ws.Range("M113").Value = -1998345
$ws.Range("N113").Value = -100006508
# row 133
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
# row 137
$ws.Range("H137").Value = 1500
$ws.Range("I137").Value = 1500
$ws.Range("K137").Value = 4500
$ws.Range("M137").Value = -1950
# row 138
$ws.Range("H138").Value = 8253.458000000001
$ws.Range("J138").Value = 8253.458000000001
$ws.Range("L138").Value = 24760.374
$ws.Range("N138").Value = -35040.374

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1280.25
$ws.Range("I2").Value = 1305.5
$ws.Range("K2").Value = 1305.5
$ws.Range("M2").Value = -1192.5
# row 32
$ws.Range("H32").Value = 3231.1538
$ws.Range("I32").Value = 3409.6365
$ws.Range("J32").Value = 2249.5
$ws.Range("K32").Value = 3409.6365
$ws.Range("L32").Value = 2249.5
$ws.Range("M32").Value = -3122.6365
$ws.Range("N32").Value = -2823.5
# row 45
$ws.Range("H45").Value = 2207
$ws.Range("J45").Value = 1269
$ws.Range("L45").Value = 1269
$ws.Range("N45").Value = -2023
# row 88
$ws.Range("H88").Value = 4666.3335
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 4499.5
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 4499.5
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -5311.5
# row 91
$ws.Range("H91").Value = 4666.3335
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 4499.5
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 4499.5
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -7307.5
# row 97
$ws.Range("H97").Value = 983.8
$ws.Range("I97").Value = 359.5
$ws.Range("K97").Value = 359.5
$ws.Range("M97").Value = 136.5
# row 102
$ws.Range("H102").Value = 2924.25
$ws.Range("I102").Value = 2924.25
$ws.Range("K102").Value = 2924.25
$ws.Range("M102").Value = -1302.25
# row 110
$ws.Range("H110").Value = 629.8
$ws.Range("I110").Value = 629.8
$ws.Range("K110").Value = 629.8
$ws.Range("M110").Value = 1415.2
# row 116
$ws.Range("H116").Value = 1280.25
$ws.Range("I116").Value = 1305.5
$ws.Range("K116").Value = 1305.5
$ws.Range("M116").Value = 988.5
# row 122
$ws.Range("H122").Value = 7265.25
$ws.Range("I122").Value = 7265.25
$ws.Range("K122").Value = 21795.75
$ws.Range("M122").Value = -19345.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1280.25
$ws.Range("I3").Value = 1305.5
$ws.Range("K3").Value = 1305.5
$ws.Range("M3").Value = -1191.5
# row 94
$ws.Range("H94").Value = 543.75
$ws.Range("I94").Value = 475
$ws.Range("K94").Value = 475
$ws.Range("M94").Value = -24
# row 105
$ws.Range("H105").Value = 1991
$ws.Range("I105").Value = 1991
$ws.Range("K105").Value = 1991
$ws.Range("M105").Value = -244
# row 134
$ws.Range("H134").Value = 4106.643
$ws.Range("I134").Value = 4076.3845
$ws.Range("K134").Value = 12229.1535
$ws.Range("M134").Value = -9694.1535

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1141.625
$ws.Range("I31").Value = 1169.1428
$ws.Range("K31").Value = 1169.1428
$ws.Range("M31").Value = -874.1428000000001
# row 34
$ws.Range("H34").Value = 1141.625
$ws.Range("I34").Value = 1169.1428
$ws.Range("K34").Value = 1169.1428
$ws.Range("M34").Value = -967.1428000000001
# row 99
$ws.Range("H99").Value = 6299.3335
$ws.Range("I99").Value = 6949
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 6949
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -5451
$ws.Range("N99").Value = -7996
# row 107
$ws.Range("H107").Value = 1738.2
$ws.Range("I107").Value = 1422.75
$ws.Range("K107").Value = 1422.75
$ws.Range("M107").Value = 497.25
# row 122
$ws.Range("H122").Value = 2136.5
$ws.Range("J122").Value = 2478.5
$ws.Range("L122").Value = 7435.5
$ws.Range("N122").Value = -12335.5
# row 126
$ws.Range("H126").Value = 6299.3335
$ws.Range("I126").Value = 6949
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 20847
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -18377
$ws.Range("N126").Value = -19940

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 132
$ws.Range("H132").Value = 2416.6667
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 2375
$ws.Range("K132").Value = 22500
$ws.Range("L132").Value = 21375
$ws.Range("M132").Value = -19970
$ws.Range("N132").Value = -26435

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3265.3333
$ws.Range("I102").Value = 3265.3333
$ws.Range("K102").Value = 3265.3333
$ws.Range("M102").Value = -1643.3333
# row 122
$ws.Range("H122").Value = 1127
$ws.Range("I122").Value = 996.6667
$ws.Range("J122").Value = 1224.75
$ws.Range("K122").Value = 2990.0001
$ws.Range("L122").Value = 3674.25
$ws.Range("M122").Value = -540.0001000000002
$ws.Range("N122").Value = -8574.25
# row 126
$ws.Range("H126").Value = 4799.4
$ws.Range("I126").Value = 3833
$ws.Range("J126").Value = 6249
$ws.Range("K126").Value = 11499
$ws.Range("L126").Value = 18747
$ws.Range("M126").Value = -9029
$ws.Range("N126").Value = -23687
# row 132
$ws.Range("H132").Value = 2408.111
$ws.Range("I132").Value = 1279.1666
$ws.Range("K132").Value = 3837.4998
$ws.Range("M132").Value = -1307.4998

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# row 40
$ws.Range("H40").Value = 4500
$ws.Range("I40").Value = 4500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4364
$ws.Range("N40").ClearContents()
# row 122
$ws.Range("H122").Value = 3497.5
$ws.Range("I122").Value = 3497.5
$ws.Range("K122").Value = 10492.5
$ws.Range("M122").Value = -8042.5
# row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 280.33334
$ws.Range("J113").Value = 232
$ws.Range("L113").Value = 696
$ws.Range("N113").Value = -5036
# row 122
$ws.Range("H122").Value = 1999.5
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("M122").Value = -3548.5
# row 126
$ws.Range("H126").Value = 1528
$ws.Range("I126").Value = 449
$ws.Range("K126").Value = 1347
$ws.Range("M126").Value = 1123
